# Applies the "Removed local data files & cleaned script from RHO files"
# edit to the DataIntoR_RHO document:
#   1. The hard-coded local setwd() path is replaced with a generic one.
#   2. The manual line break right after read.csv("Iris.csv") is turned
#      into a real paragraph break, so "> str(dfobj)" becomes its own
#      SourceCode paragraph instead of being crammed onto the same line.

$d = $word.ActiveDocument

# --- 1. Swap the local working-directory literal in setwd(...) ----------
$oldPath = '"C:/aaaWork/Web/GitHub/NCMTH107/modules/"'
$newPath = '"C:/stats/"'

$full = $d.Content.Text
$idx = $full.IndexOf($oldPath)
if ($idx -ge 0) {
    # Set .Text directly (rather than Find/Replace) so Word's smart-quote
    # autocorrect doesn't mangle the straight quotes inside the literal.
    $r = $d.Range($idx, $idx + $oldPath.Length)
    $r.Text = $newPath
}

# --- 2. Turn the "Iris.csv")<br>" manual break into a paragraph break ---
$marker = 'Iris.csv")'
$full = $d.Content.Text
$midx = $full.IndexOf($marker)
if ($midx -ge 0) {
    $breakPos = $midx + $marker.Length
    # The character immediately after ")" is the manual line-break run;
    # replacing it with a carriage return splits the paragraph in two and
    # the new paragraph inherits the SourceCode style automatically.
    $br = $d.Range($breakPos, $breakPos + 1)
    $br.Text = "`r"
}
